$p = $ppt.ActivePresentation

# Remove the last two slides (the "Views" slide and the "USERSCONTROLLER.php" slide)
$count = $p.Slides.Count
$p.Slides.Item($count).Delete()
$p.Slides.Item($count - 1).Delete()
